$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "31/12/2006"
$ws.Range("C2").Value = 1.87984550086694

$ws.Range("B3").Value = "31/12/2010"
$ws.Range("C3").Value = 1.74571102459612

$ws.Range("B4").Value = "31/12/2014"
$ws.Range("C4").Value = 1.78994618138214

$ws.Range("B5").Value = "31/12/2018"
$ws.Range("C5").Value = 1.77923455058482

$ws.Range("B6").Value = "31/12/2022"
$ws.Range("C6").Value = 1.57625383660638

$ws.Range("B7").Value = "31/12/2024"
$ws.Range("C7").Value = 1.54680260681053

$ws.Range("B8").Value = "31/12/2006"
$ws.Range("C8").Value = 2.06230151318107

$ws.Range("B9").Value = "31/12/2010"
$ws.Range("C9").Value = 1.83798130257565

$ws.Range("B10").Value = "31/12/2014"
$ws.Range("C10").Value = 1.77809089870609

$ws.Range("B11").Value = "31/12/2018"
$ws.Range("C11").Value = 1.80512077932592

$ws.Range("B12").Value = "31/12/2022"
$ws.Range("C12").Value = 1.56612675297142

$ws.Range("B13").Value = "31/12/2024"
$ws.Range("C13").Value = 1.54418474547815

$ws.Range("B14").Value = "31/12/2006"
$ws.Range("C14").Value = 2.08844892856124

$ws.Range("B15").Value = "31/12/2010"
$ws.Range("C15").Value = 1.8277327274294

$ws.Range("B16").Value = "31/12/2014"
$ws.Range("C16").Value = 1.78769503823641

$ws.Range("B17").Value = "31/12/2018"
$ws.Range("C17").Value = 1.78826541448712

$ws.Range("B18").Value = "31/12/2022"
$ws.Range("C18").Value = 1.52533730886118

$ws.Range("B19").Value = "31/12/2024"
$ws.Range("C19").Value = 1.52462190148086
